# Update the quarterly income-statement database: drop the oldest quarter
# (column D, "فصل دوم منتهی به 1399/06") and append the newest quarter
# ("فصل چهارم منتهی به 1401/12") as a new column M, shifting everything
# else one column to the left. Also correct the revised publish date for
# the Q3-1401 report and one data-entry fix in the "سرمایه" row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Drop the oldest quarter column (D). Everything to its right (E:M)
#    shifts left to become the new D:L, carrying values/styles with it.
$ws.Columns.Item(4).Delete()

# 2) Bring in a brand-new column M with the same look as the rest of the
#    table (copy number formats/styles from the now-last data column L),
#    and give it the wider "latest quarter" column width (raw width 31).
$ws.Range("L1:L28").Copy()
$ws.Range("M1:M28").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Columns.Item(13).ColumnWidth = 30.1666667

# 3) Populate the new quarter's header/publish-date/data in column M.
$ws.Range("M8").Value = "فصل چهارم منتهی به 1401/12"
$ws.Range("M9").Value = "1402-02-30"

$ws.Range("M11").Value = 5626489
$ws.Range("M12").Value = -5520040
$ws.Range("M13").Value = 106449
$ws.Range("M14").Value = -91747
$ws.Range("M15").Value = 0
$ws.Range("M16").Value = 3216
$ws.Range("M17").Value = 17918
$ws.Range("M18").Value = -10500
$ws.Range("M19").Value = 137150
$ws.Range("M20").Value = 144568
$ws.Range("M21").Value = 61440
$ws.Range("M22").Value = 206008
$ws.Range("M23").Value = 0
$ws.Range("M24").Value = 206008
$ws.Range("M25").Value = 1058
$ws.Range("M26").Value = 194650
$ws.Range("M27").Value = 1056

# 4) The Q3-1401 report's publish date was revised/updated (now column I).
$ws.Range("I9").Value = "1402-02-30 (8)"

# 5) Data correction in the "سرمایه" (capital) row for the Q2-1401 column
#    (now column I) picked up during this refresh.
$ws.Range("I26").Value = 194650
